$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# --- Copy number-format styling from row 43 (D:F only - avoids touching
#     column A/C/G/H, which keeps the new rows' blank cells truly empty) ---
$ws.Range("D43:F43").Copy()
$ws.Range("D44:F46").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- New rows of "default" baseline test runs ---
$ws.Range("A44").Value = "default.1"
$ws.Range("A45").Value = "default.2"
$ws.Range("A46").Value = "default.3"
$ws.Range("B44:B46").Value = "default"

$ws.Range("D44:D46").Formula = "=1/3"
$ws.Range("E44").Formula = "=C44/D44"
$ws.Range("E45").Formula = "=C45/D45"
$ws.Range("E46").Formula = "=C46/D46"
$ws.Range("F44").Formula = "=AVERAGEIF(B:B,B44,E:E)"
$ws.Range("F45").Formula = "=AVERAGEIF(B:B,B45,E:E)"
$ws.Range("F46").Formula = "=AVERAGEIF(B:B,B46,E:E)"

# --- Note explaining the new "default" group ---
$comment = $ws.Range("B44").AddComment("default:" + [char]10 + "Create a default set of weights and a model for testing to see which songs are mis-identified every time.")

# --- Selection left where the author's cursor ended up ---
$ws.Range("G43").Select()

Write-Host "done"
